$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) column cells being updated so that
# numeric-looking strings (e.g. "2.470", "1.000") keep their exact
# text representation (incl. trailing zeros / dot-grouping), matching
# the original inline-string cell contents, instead of Excel silently
# coercing the assignment into a Number.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "23.452.82"
$ws.Range("D3").Value = "1.629.75"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").Value = "304.39"
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("D7").Value = "0.3769"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "0.3659"
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").Value = "51.83"
$ws.Range("E9").Value = "  -1.34%  "
$ws.Range("D10").Value = "0.08215"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").Value = "1.225"
$ws.Range("E11").Value = "  -3.58%  "
$ws.Range("D13").Value = "22.44"
$ws.Range("E13").Value = "  -2.55%  "
$ws.Range("D14").Value = "6.557"
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("D15").Value = "0.00001252"
$ws.Range("E15").Value = "  -2.11%  "
$ws.Range("D16").Value = "7.257"
$ws.Range("E16").Value = "  -1.95%  "
$ws.Range("D17").Value = "1.628.72"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").Value = "94.12"
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("E20").Value = "  -2.90%  "
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "12.71"
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("D24").Value = "23.450.29"
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D25").Value = "3.202"
$ws.Range("E25").Value = "  +3.36%  "
$ws.Range("D26").Value = "2.470"
$ws.Range("E26").Value = "  +1.88%  "
$ws.Range("E27").Value = "  +0.40%  "
$ws.Range("D28").Value = "149.95"
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("D29").Value = "5.312"
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("D30").Value = "133.88"
$ws.Range("E30").Value = "  -1.26%  "
$ws.Range("D31").Value = "1.810.41"
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("D32").Value = "2.275"
$ws.Range("E32").Value = "  -3.93%  "
$ws.Range("D33").Value = "6.809"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").Value = "1.017"
$ws.Range("E34").Value = "  +5.31%  "
$ws.Range("D35").Value = "10.79"
$ws.Range("E35").Value = "  +4.36%  "
$ws.Range("D36").Value = "0.02787"
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("D38").Value = "0.08749"
$ws.Range("E38").Value = "  -1.36%  "
$ws.Range("D39").Value = "0.07127"
$ws.Range("E39").Value = "  -3.23%  "
$ws.Range("D40").Value = "6.036"
$ws.Range("E40").Value = "  -2.50%  "
$ws.Range("D41").Value = "0.7054"
$ws.Range("E41").Value = "  -0.83%  "
$ws.Range("E42").Value = "  -2.32%  "
$ws.Range("D43").Value = "16.28"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("E44").Value = "  -2.11%  "
$ws.Range("D45").Value = "0.6556"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("E46").Value = "  -0.83%  "
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").Value = "3.985"
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("D49").Value = "0.08019"
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("D50").Value = "1.202"
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("D51").Value = "125.86"
$ws.Range("E51").Value = "  -2.77%  "
